$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the contract rows (A2:K41) ascending by Start Date (column B),
# using the worksheet's Sort object so the sort state/condition is
# persisted to the file, matching Excel's own "Sort" UI action.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B41"), 0, 1, 0, 0)
$ws.Sort.SetRange($ws.Range("A2:K41"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Move the selection/active cell to B2 and reset the scroll position so
# the view no longer shows the old "topLeftCell=L1 / L1:AF1048576"
# selection left over from before the sort.
$ws.Range("B2").Select()
